# The workbook's active sheet is "이정원" (the 4th tab, sheetId=4 /
# Worksheets.Item(4)) — this is the sheet the diff touches.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Row 4 of the "할일/내용/할당일/완성일/결과/문제점" task table was blank
# (styled placeholder cells only). Fill it in with a new task entry:
#   할일(A)  = "기능 수정"
#   내용(B)  = "결제창 생성, 회원 검색 기능과 현금 결제, 카드 결제 구분 기능 추가"
#   할당일(C) = 2019-05-10 (serial 43595)
#   완성일(D) = 2019-05-11 (serial 43596)
#   결과(E)/문제점(F) reuse the same text already used in row 3, so read
#   those live (via the method-call form of .Value so we get the actual
#   string back, not a property descriptor) instead of retyping them —
#   this keeps them pointing at the existing shared-string entries.
$ws.Range("A4").Value = "기능 수정"
$ws.Range("B4").Value = "결제창 생성, 회원 검색 기능과 현금 결제, 카드 결제 구분 기능 추가"
$ws.Range("C4").Value = 43595
$ws.Range("D4").Value = 43596
$ws.Range("E4").Value = $ws.Range("E3").Value()
$ws.Range("F4").Value = $ws.Range("F3").Value()

# The new row's content wraps onto multiple lines, so the row is taller
# than the default — set its height to match.
$ws.Rows.Item(4).RowHeight = 70
